$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195 (shifts existing rows 195-226 down to 196-227)
$ws.Rows("195:195").Insert()

# Populate the new row 195 with a new weekly price observation
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44505
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 100112043
$ws.Range("G195").Value = "Pepino ensalada"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 95
$ws.Range("K195").Value = 7000
$ws.Range("L195").Value = 7500
$ws.Range("M195").Value = 7237
$ws.Range("N195").Value = "`$/caja 70 unidades"
$ws.Range("O195").Value = "Región de Arica y Parinacota"
$ws.Range("P195").Value = 103
$ws.Range("Q195").Value = 70
$ws.Range("R195").Value = "Hortaliza"
